$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Sector Cap" row (old row 8); everything below shifts up by one row.
$ws.Rows.Item(8).Delete()

# Relabel the performance rows and refresh every sector value (new calc run).
$ws.Range("A2").Value = "Pct. Return 1M"
$ws.Range("B2").Value = -0.169
$ws.Range("C2").Value = 3.406
$ws.Range("D2").Value = 6.72
$ws.Range("E2").Value = -7.431
$ws.Range("F2").Value = 5.361
$ws.Range("G2").Value = 7.388
$ws.Range("H2").Value = -7.003
$ws.Range("I2").Value = -8.436999999999999
$ws.Range("J2").Value = 11.294
$ws.Range("K2").Value = 13.503
$ws.Range("L2").Value = 1.849

$ws.Range("A3").Value = "Pct. Return 3M"
$ws.Range("B3").Value = -8.271000000000001
$ws.Range("C3").Value = -12.956
$ws.Range("D3").Value = -30.77
$ws.Range("E3").Value = -26.267
$ws.Range("F3").Value = -19.926
$ws.Range("G3").Value = -12.543
$ws.Range("H3").Value = -46.571
$ws.Range("I3").Value = -34.373
$ws.Range("J3").Value = -1.492
$ws.Range("K3").Value = -3.074
$ws.Range("L3").Value = -18.034

$ws.Range("A4").Value = "Pct. Return YTD"
$ws.Range("B4").Value = -9.782
$ws.Range("C4").Value = -8.561
$ws.Range("D4").Value = -28.539
$ws.Range("E4").Value = -18.944
$ws.Range("F4").Value = -15.487
$ws.Range("G4").Value = -7.514
$ws.Range("H4").Value = -43.338
$ws.Range("I4").Value = -28.168
$ws.Range("J4").Value = -3.654
$ws.Range("K4").Value = 6.349
$ws.Range("L4").Value = -14.58

$ws.Range("A5").Value = "Pct. Return 1Y"
$ws.Range("B5").Value = -8.156000000000001
$ws.Range("C5").Value = 1.379
$ws.Range("D5").Value = -30.691
$ws.Range("E5").Value = -8.654
$ws.Range("F5").Value = -11.97
$ws.Range("G5").Value = 1.614
$ws.Range("H5").Value = -35.762
$ws.Range("I5").Value = -18.721
$ws.Range("J5").Value = 6.867
$ws.Range("K5").Value = 18.768
$ws.Range("L5").Value = -16.698

$ws.Range("A6").Value = "Pct. Return Max"
$ws.Range("B6").Value = -1.446
$ws.Range("C6").Value = 46.582
$ws.Range("D6").Value = -34.693
$ws.Range("E6").Value = 35.836
$ws.Range("F6").Value = 37.168
$ws.Range("G6").Value = 19.902
$ws.Range("H6").Value = -40.495
$ws.Range("I6").Value = -7.879
$ws.Range("J6").Value = 18.361
$ws.Range("K6").Value = 94.608
$ws.Range("L6").Value = -20.371

$ws.Range("A7").Value = "Pct. Return CAGR"
$ws.Range("B7").Value = -0.314
$ws.Range("C7").Value = 8.621
$ws.Range("D7").Value = -8.801
$ws.Range("E7").Value = 6.847
$ws.Range("F7").Value = 7.073
$ws.Range("G7").Value = 4.003
$ws.Range("H7").Value = -10.618
$ws.Range("I7").Value = -1.759
$ws.Range("J7").Value = 3.712
$ws.Range("K7").Value = 15.485
$ws.Range("L7").Value = -4.806

$ws.Range("A8").Value = "Beta"
$ws.Range("B8").Value = 0.756
$ws.Range("C8").Value = 0.985
$ws.Range("D8").Value = 1.044
$ws.Range("E8").Value = 0.979
$ws.Range("F8").Value = 1.009
$ws.Range("G8").Value = 1.035
$ws.Range("H8").Value = 1.115
$ws.Range("I8").Value = 0.859
$ws.Range("J8").Value = 0.612
$ws.Range("K8").Value = 1.039
$ws.Range("L8").Value = 0.604

$ws.Range("A9").Value = "Alpha"
$ws.Range("B9").Value = -6.515
$ws.Range("C9").Value = 0.734
$ws.Range("D9").Value = -17.124
$ws.Range("E9").Value = -1.002
$ws.Range("F9").Value = -0.996
$ws.Range("G9").Value = -4.255
$ws.Range("H9").Value = -19.462
$ws.Range("I9").Value = -8.720000000000001
$ws.Range("J9").Value = -1.435
$ws.Range("K9").Value = 7.196
$ws.Range("L9").Value = -9.891

$ws.Range("A10").Value = "Sharpe"
$ws.Range("B10").Value = -0.054
$ws.Range("C10").Value = 0.394
$ws.Range("D10").Value = -0.291
$ws.Range("E10").Value = 0.255
$ws.Range("F10").Value = 0.292
$ws.Range("G10").Value = 0.138
$ws.Range("H10").Value = -0.412
$ws.Range("I10").Value = -0.097
$ws.Range("J10").Value = 0.269
$ws.Range("K10").Value = 0.621
$ws.Range("L10").Value = -0.428

$ws.Range("A11").Value = "Treynor"
$ws.Range("B11").Value = -0.013
$ws.Range("C11").Value = 0.081
$ws.Range("D11").Value = -0.09
$ws.Range("E11").Value = 0.063
$ws.Range("F11").Value = 0.064
$ws.Range("G11").Value = 0.032
$ws.Range("H11").Value = -0.101
$ws.Range("I11").Value = -0.028
$ws.Range("J11").Value = 0.05
$ws.Range("K11").Value = 0.143
$ws.Range("L11").Value = -0.09

$ws.Range("A12").Value = "Max Drawdown"
$ws.Range("B12").Value = -39.365
$ws.Range("C12").Value = -32.427
$ws.Range("D12").Value = -66.36499999999999
$ws.Range("E12").Value = -44.101
$ws.Range("F12").Value = -39.541
$ws.Range("G12").Value = -36.075
$ws.Range("H12").Value = -50.385
$ws.Range("I12").Value = -42.377
$ws.Range("J12").Value = -22.737
$ws.Range("K12").Value = -34.126
$ws.Range("L12").Value = -34.779

$ws.Range("A13").Value = "Std. Deviation"
$ws.Range("B13").Value = 17.826
$ws.Range("C13").Value = 20.25
$ws.Range("D13").Value = 32.513
$ws.Range("E13").Value = 24.33
$ws.Range("F13").Value = 22.009
$ws.Range("G13").Value = 24.327
$ws.Range("H13").Value = 27.329
$ws.Range("I13").Value = 24.775
$ws.Range("J13").Value = 11.445
$ws.Range("K13").Value = 23.917
$ws.Range("L13").Value = 12.739

$ws.Range("A14").Value = "R Squared"
$ws.Range("B14").Value = 0.344
$ws.Range("C14").Value = 0.761
$ws.Range("D14").Value = 0.471
$ws.Range("E14").Value = 0.416
$ws.Range("F14").Value = 0.791
$ws.Range("G14").Value = 0.519
$ws.Range("H14").Value = 0.471
$ws.Range("I14").Value = 0.32
$ws.Range("J14").Value = 0.375
$ws.Range("K14").Value = 0.742
$ws.Range("L14").Value = 0.288

$ws.Range("A15").Value = "Expected Return"
$ws.Range("B15").Value = 6.201
$ws.Range("C15").Value = 7.887
$ws.Range("D15").Value = 8.321999999999999
$ws.Range("E15").Value = 7.849
$ws.Range("F15").Value = 8.069000000000001
$ws.Range("G15").Value = 8.257999999999999
$ws.Range("H15").Value = 8.843999999999999
$ws.Range("I15").Value = 6.961
$ws.Range("J15").Value = 5.147
$ws.Range("K15").Value = 8.289
$ws.Range("L15").Value = 5.085
